$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 330, shifting the existing rows (330-416) down to (332-418).
# Excel automatically copies formatting (e.g. the date NumberFormat on column D) from the row above.
$ws.Rows("330:331").Insert()

# Fill in the two new rows (330 and 331) with the new record data.
# Columns A,B,C,E,F,G,H,I,R keep the same values as the surrounding rows for this market/product.

# New row 330
$ws.Range("A330").Value = 3
$ws.Range("B330").Value = "Femacal de La Calera"
$ws.Range("C330").Value = "Coquimbo"
$ws.Range("D330").Value2 = 44642
$ws.Range("E330").Value = 5
$ws.Range("F330").Value = 100112032
$ws.Range("G330").Value = "Zapallo italiano"
$ws.Range("H330").Value = "Sin especificar"
$ws.Range("I330").Value = "Primera"
$ws.Range("J330").Value = 130
$ws.Range("K330").Value = 4000
$ws.Range("L330").Value = 4500
$ws.Range("M330").Value = 4269
$ws.Range("N330").Value = "$/caja 36 unidades"
$ws.Range("O330").Value = "Provincia de Quillota"
$ws.Range("P330").Value = 119
$ws.Range("Q330").Value = 36
$ws.Range("R330").Value = "Hortaliza"

# New row 331
$ws.Range("A331").Value = 3
$ws.Range("B331").Value = "Femacal de La Calera"
$ws.Range("C331").Value = "Coquimbo"
$ws.Range("D331").Value2 = 44642
$ws.Range("E331").Value = 5
$ws.Range("F331").Value = 100112032
$ws.Range("G331").Value = "Zapallo italiano"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 160
$ws.Range("K331").Value = 9000
$ws.Range("L331").Value = 9500
$ws.Range("M331").Value = 9156
$ws.Range("N331").Value = "$/caja 70 unidades"
$ws.Range("O331").Value = "Provincia de Quillota"
$ws.Range("P331").Value = 131
$ws.Range("Q331").Value = 70
$ws.Range("R331").Value = "Hortaliza"
